$wb = $excel.ActiveWorkbook
$mismatchCount = 0

# Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
if ($ws.Range('I2').Value2 -ne 6637) { $mismatchCount++; Write-Output "MISMATCH Citywide Totals!I2: expected 6637 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 6655
if ($ws.Range('I3').Value2 -ne 6944) { $mismatchCount++; Write-Output "MISMATCH Citywide Totals!I3: expected 6944 got $($ws.Range('I3').Value2)" }
$ws.Range('I3').Value = 6955
if ($ws.Range('I4').Value2 -ne 1589) { $mismatchCount++; Write-Output "MISMATCH Citywide Totals!I4: expected 1589 got $($ws.Range('I4').Value2)" }
$ws.Range('I4').Value = 1591
if ($ws.Range('I5').Value2 -ne 645) { $mismatchCount++; Write-Output "MISMATCH Citywide Totals!I5: expected 645 got $($ws.Range('I5').Value2)" }
$ws.Range('I5').Value = 648
if ($ws.Range('I6').Value2 -ne 7992) { $mismatchCount++; Write-Output "MISMATCH Citywide Totals!I6: expected 7992 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 8028
if ($ws.Range('I7').Value2 -ne 23807) { $mismatchCount++; Write-Output "MISMATCH Citywide Totals!I7: expected 23807 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 23877

# Uptown
$ws = $wb.Worksheets.Item('Uptown')
if ($ws.Range('I6').Value2 -ne 98) { $mismatchCount++; Write-Output "MISMATCH Uptown!I6: expected 98 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 99
if ($ws.Range('I7').Value2 -ne 283) { $mismatchCount++; Write-Output "MISMATCH Uptown!I7: expected 283 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 284

# Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
if ($ws.Range('I6').Value2 -ne 223) { $mismatchCount++; Write-Output "MISMATCH Grand Crossing!I6: expected 223 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 226
if ($ws.Range('I7').Value2 -ne 748) { $mismatchCount++; Write-Output "MISMATCH Grand Crossing!I7: expected 748 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 751

# South Deering
$ws = $wb.Worksheets.Item('South Deering')
if ($ws.Range('I3').Value2 -ne 70) { $mismatchCount++; Write-Output "MISMATCH South Deering!I3: expected 70 got $($ws.Range('I3').Value2)" }
$ws.Range('I3').Value = 71
if ($ws.Range('I7').Value2 -ne 211) { $mismatchCount++; Write-Output "MISMATCH South Deering!I7: expected 211 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 212

# New City
$ws = $wb.Worksheets.Item('New City')
if ($ws.Range('I3').Value2 -ne 162) { $mismatchCount++; Write-Output "MISMATCH New City!I3: expected 162 got $($ws.Range('I3').Value2)" }
$ws.Range('I3').Value = 163
if ($ws.Range('I7').Value2 -ne 549) { $mismatchCount++; Write-Output "MISMATCH New City!I7: expected 549 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 550

# By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
if ($ws.Range('I7').Value2 -ne 749) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I7: expected 749 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 750
if ($ws.Range('I8').Value2 -ne 1430) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I8: expected 1430 got $($ws.Range('I8').Value2)" }
$ws.Range('I8').Value = 1432
if ($ws.Range('I10').Value2 -ne 170) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I10: expected 170 got $($ws.Range('I10').Value2)" }
$ws.Range('I10').Value = 171
if ($ws.Range('I11').Value2 -ne 361) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I11: expected 361 got $($ws.Range('I11').Value2)" }
$ws.Range('I11').Value = 362
if ($ws.Range('I12').Value2 -ne 57) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I12: expected 57 got $($ws.Range('I12').Value2)" }
$ws.Range('I12').Value = 58
if ($ws.Range('I15').Value2 -ne 276) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I15: expected 276 got $($ws.Range('I15').Value2)" }
$ws.Range('I15').Value = 277
if ($ws.Range('I18').Value2 -ne 182) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I18: expected 182 got $($ws.Range('I18').Value2)" }
$ws.Range('I18').Value = 184
if ($ws.Range('I19').Value2 -ne 671) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I19: expected 671 got $($ws.Range('I19').Value2)" }
$ws.Range('I19').Value = 673
if ($ws.Range('I20').Value2 -ne 590) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I20: expected 590 got $($ws.Range('I20').Value2)" }
$ws.Range('I20').Value = 591
if ($ws.Range('I25').Value2 -ne 124) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I25: expected 124 got $($ws.Range('I25').Value2)" }
$ws.Range('I25').Value = 125
if ($ws.Range('I29').Value2 -ne 1437) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I29: expected 1437 got $($ws.Range('I29').Value2)" }
$ws.Range('I29').Value = 1440
if ($ws.Range('I33').Value2 -ne 1065) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I33: expected 1065 got $($ws.Range('I33').Value2)" }
$ws.Range('I33').Value = 1066
if ($ws.Range('I36').Value2 -ne 325) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I36: expected 325 got $($ws.Range('I36').Value2)" }
$ws.Range('I36').Value = 326
if ($ws.Range('I37').Value2 -ne 748) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I37: expected 748 got $($ws.Range('I37').Value2)" }
$ws.Range('I37').Value = 751
if ($ws.Range('I42').Value2 -ne 873) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I42: expected 873 got $($ws.Range('I42').Value2)" }
$ws.Range('I42').Value = 881
if ($ws.Range('I44').Value2 -ne 176) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I44: expected 176 got $($ws.Range('I44').Value2)" }
$ws.Range('I44').Value = 179
if ($ws.Range('I46').Value2 -ne 53) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I46: expected 53 got $($ws.Range('I46').Value2)" }
$ws.Range('I46').Value = 54
if ($ws.Range('I49').Value2 -ne 158) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I49: expected 158 got $($ws.Range('I49').Value2)" }
$ws.Range('I49').Value = 159
if ($ws.Range('I52').Value2 -ne 529) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I52: expected 529 got $($ws.Range('I52').Value2)" }
$ws.Range('I52').Value = 542
if ($ws.Range('I54').Value2 -ne 480) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I54: expected 480 got $($ws.Range('I54').Value2)" }
$ws.Range('I54').Value = 481
if ($ws.Range('I60').Value2 -ne 134) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I60: expected 134 got $($ws.Range('I60').Value2)" }
$ws.Range('I60').Value = 135
if ($ws.Range('I63').Value2 -ne 71) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I63: expected 71 got $($ws.Range('I63').Value2)" }
$ws.Range('I63').Value = 73
if ($ws.Range('I65').Value2 -ne 549) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I65: expected 549 got $($ws.Range('I65').Value2)" }
$ws.Range('I65').Value = 550
if ($ws.Range('I66').Value2 -ne 65) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I66: expected 65 got $($ws.Range('I66').Value2)" }
$ws.Range('I66').Value = 66
if ($ws.Range('I68').Value2 -ne 80) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I68: expected 80 got $($ws.Range('I68').Value2)" }
$ws.Range('I68').Value = 81
if ($ws.Range('I76').Value2 -ne 344) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I76: expected 344 got $($ws.Range('I76').Value2)" }
$ws.Range('I76').Value = 345
if ($ws.Range('I78').Value2 -ne 320) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I78: expected 320 got $($ws.Range('I78').Value2)" }
$ws.Range('I78').Value = 321
if ($ws.Range('I79').Value2 -ne 677) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I79: expected 677 got $($ws.Range('I79').Value2)" }
$ws.Range('I79').Value = 679
if ($ws.Range('I83').Value2 -ne 515) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I83: expected 515 got $($ws.Range('I83').Value2)" }
$ws.Range('I83').Value = 518
if ($ws.Range('I84').Value2 -ne 211) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I84: expected 211 got $($ws.Range('I84').Value2)" }
$ws.Range('I84').Value = 212
if ($ws.Range('I85').Value2 -ne 1067) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I85: expected 1067 got $($ws.Range('I85').Value2)" }
$ws.Range('I85').Value = 1070
if ($ws.Range('I86').Value2 -ne 150) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I86: expected 150 got $($ws.Range('I86').Value2)" }
$ws.Range('I86').Value = 151
if ($ws.Range('I88').Value2 -ne 221) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I88: expected 221 got $($ws.Range('I88').Value2)" }
$ws.Range('I88').Value = 222
if ($ws.Range('I89').Value2 -ne 283) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I89: expected 283 got $($ws.Range('I89').Value2)" }
$ws.Range('I89').Value = 284
if ($ws.Range('I95').Value2 -ne 361) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I95: expected 361 got $($ws.Range('I95').Value2)" }
$ws.Range('I95').Value = 362
if ($ws.Range('I97').Value2 -ne 194) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I97: expected 194 got $($ws.Range('I97').Value2)" }
$ws.Range('I97').Value = 195
if ($ws.Range('I101').Value2 -ne 23807) { $mismatchCount++; Write-Output "MISMATCH By Neighborhood!I101: expected 23807 got $($ws.Range('I101').Value2)" }
$ws.Range('I101').Value = 23877

# South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
if ($ws.Range('I2').Value2 -ne 174) { $mismatchCount++; Write-Output "MISMATCH South Chicago!I2: expected 174 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 175
if ($ws.Range('I3').Value2 -ne 185) { $mismatchCount++; Write-Output "MISMATCH South Chicago!I3: expected 185 got $($ws.Range('I3').Value2)" }
$ws.Range('I3').Value = 186
if ($ws.Range('I6').Value2 -ne 114) { $mismatchCount++; Write-Output "MISMATCH South Chicago!I6: expected 114 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 115
if ($ws.Range('I7').Value2 -ne 515) { $mismatchCount++; Write-Output "MISMATCH South Chicago!I7: expected 515 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 518

# West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
if ($ws.Range('I3').Value2 -ne 128) { $mismatchCount++; Write-Output "MISMATCH West Pullman!I3: expected 128 got $($ws.Range('I3').Value2)" }
$ws.Range('I3').Value = 129
if ($ws.Range('I7').Value2 -ne 361) { $mismatchCount++; Write-Output "MISMATCH West Pullman!I7: expected 361 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 362

# Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
if ($ws.Range('I2').Value2 -ne 239) { $mismatchCount++; Write-Output "MISMATCH Garfield Park!I2: expected 239 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 240
if ($ws.Range('I7').Value2 -ne 1065) { $mismatchCount++; Write-Output "MISMATCH Garfield Park!I7: expected 1065 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 1066

# Lincoln Park
$ws = $wb.Worksheets.Item('Lincoln Park')
if ($ws.Range('I3').Value2 -ne 16) { $mismatchCount++; Write-Output "MISMATCH Lincoln Park!I3: expected 16 got $($ws.Range('I3').Value2)" }
$ws.Range('I3').Value = 17
if ($ws.Range('I7').Value2 -ne 158) { $mismatchCount++; Write-Output "MISMATCH Lincoln Park!I7: expected 158 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 159

# Loop
$ws = $wb.Worksheets.Item('Loop')
if ($ws.Range('I3').Value2 -ne 103) { $mismatchCount++; Write-Output "MISMATCH Loop!I3: expected 103 got $($ws.Range('I3').Value2)" }
$ws.Range('I3').Value = 104
if ($ws.Range('I7').Value2 -ne 480) { $mismatchCount++; Write-Output "MISMATCH Loop!I7: expected 480 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 481

# Englewood
$ws = $wb.Worksheets.Item('Englewood')
if ($ws.Range('I6').Value2 -ne 394) { $mismatchCount++; Write-Output "MISMATCH Englewood!I6: expected 394 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 397
if ($ws.Range('I7').Value2 -ne 1437) { $mismatchCount++; Write-Output "MISMATCH Englewood!I7: expected 1437 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 1440

# Chatham
$ws = $wb.Worksheets.Item('Chatham')
if ($ws.Range('I3').Value2 -ne 195) { $mismatchCount++; Write-Output "MISMATCH Chatham!I3: expected 195 got $($ws.Range('I3').Value2)" }
$ws.Range('I3').Value = 196
if ($ws.Range('I5').Value2 -ne 18) { $mismatchCount++; Write-Output "MISMATCH Chatham!I5: expected 18 got $($ws.Range('I5').Value2)" }
$ws.Range('I5').Value = 19
if ($ws.Range('I7').Value2 -ne 671) { $mismatchCount++; Write-Output "MISMATCH Chatham!I7: expected 671 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 673

# Irving Park
$ws = $wb.Worksheets.Item('Irving Park')
if ($ws.Range('I2').Value2 -ne 55) { $mismatchCount++; Write-Output "MISMATCH Irving Park!I2: expected 55 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 57
if ($ws.Range('I6').Value2 -ne 50) { $mismatchCount++; Write-Output "MISMATCH Irving Park!I6: expected 50 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 51
if ($ws.Range('I7').Value2 -ne 176) { $mismatchCount++; Write-Output "MISMATCH Irving Park!I7: expected 176 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 179

# River North
$ws = $wb.Worksheets.Item('River North')
if ($ws.Range('I2').Value2 -ne 69) { $mismatchCount++; Write-Output "MISMATCH River North!I2: expected 69 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 70
if ($ws.Range('I7').Value2 -ne 344) { $mismatchCount++; Write-Output "MISMATCH River North!I7: expected 344 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 345

# South Shore
$ws = $wb.Worksheets.Item('South Shore')
if ($ws.Range('I2').Value2 -ne 301) { $mismatchCount++; Write-Output "MISMATCH South Shore!I2: expected 301 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 303
if ($ws.Range('I6').Value2 -ne 278) { $mismatchCount++; Write-Output "MISMATCH South Shore!I6: expected 278 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 279
if ($ws.Range('I7').Value2 -ne 1067) { $mismatchCount++; Write-Output "MISMATCH South Shore!I7: expected 1067 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 1070

# Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
if ($ws.Range('I2').Value2 -ne 206) { $mismatchCount++; Write-Output "MISMATCH Humboldt Park!I2: expected 206 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 207
if ($ws.Range('I3').Value2 -ne 265) { $mismatchCount++; Write-Output "MISMATCH Humboldt Park!I3: expected 265 got $($ws.Range('I3').Value2)" }
$ws.Range('I3').Value = 266
if ($ws.Range('I6').Value2 -ne 318) { $mismatchCount++; Write-Output "MISMATCH Humboldt Park!I6: expected 318 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 324
if ($ws.Range('I7').Value2 -ne 873) { $mismatchCount++; Write-Output "MISMATCH Humboldt Park!I7: expected 873 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 881

# Avondale
$ws = $wb.Worksheets.Item('Avondale')
if ($ws.Range('I6').Value2 -ne 77) { $mismatchCount++; Write-Output "MISMATCH Avondale!I6: expected 77 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 78
if ($ws.Range('I7').Value2 -ne 170) { $mismatchCount++; Write-Output "MISMATCH Avondale!I7: expected 170 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 171

# Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
if ($ws.Range('I5').Value2 -ne 7) { $mismatchCount++; Write-Output "MISMATCH Rogers Park!I5: expected 7 got $($ws.Range('I5').Value2)" }
$ws.Range('I5').Value = 8
if ($ws.Range('I7').Value2 -ne 320) { $mismatchCount++; Write-Output "MISMATCH Rogers Park!I7: expected 320 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 321

# Jefferson Park
$ws = $wb.Worksheets.Item('Jefferson Park')
if ($ws.Range('I2').Value2 -ne 14) { $mismatchCount++; Write-Output "MISMATCH Jefferson Park!I2: expected 14 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 15
if ($ws.Range('I7').Value2 -ne 53) { $mismatchCount++; Write-Output "MISMATCH Jefferson Park!I7: expected 53 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 54

# Roseland
$ws = $wb.Worksheets.Item('Roseland')
if ($ws.Range('I2').Value2 -ne 196) { $mismatchCount++; Write-Output "MISMATCH Roseland!I2: expected 196 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 197
if ($ws.Range('I4').Value2 -ne 39) { $mismatchCount++; Write-Output "MISMATCH Roseland!I4: expected 39 got $($ws.Range('I4').Value2)" }
$ws.Range('I4').Value = 40
if ($ws.Range('I7').Value2 -ne 677) { $mismatchCount++; Write-Output "MISMATCH Roseland!I7: expected 677 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 679

# Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
if ($ws.Range('I3').Value2 -ne 167) { $mismatchCount++; Write-Output "MISMATCH Chicago Lawn!I3: expected 167 got $($ws.Range('I3').Value2)" }
$ws.Range('I3').Value = 168
if ($ws.Range('I7').Value2 -ne 590) { $mismatchCount++; Write-Output "MISMATCH Chicago Lawn!I7: expected 590 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 591

# Calumet Heights
$ws = $wb.Worksheets.Item('Calumet Heights')
if ($ws.Range('I2').Value2 -ne 49) { $mismatchCount++; Write-Output "MISMATCH Calumet Heights!I2: expected 49 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 51
if ($ws.Range('I7').Value2 -ne 182) { $mismatchCount++; Write-Output "MISMATCH Calumet Heights!I7: expected 182 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 184

# Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
if ($ws.Range('I2').Value2 -ne 92) { $mismatchCount++; Write-Output "MISMATCH Grand Boulevard!I2: expected 92 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 93
if ($ws.Range('I7').Value2 -ne 325) { $mismatchCount++; Write-Output "MISMATCH Grand Boulevard!I7: expected 325 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 326

# Little Village
$ws = $wb.Worksheets.Item('Little Village')
if ($ws.Range('I6').Value2 -ne 160) { $mismatchCount++; Write-Output "MISMATCH Little Village!I6: expected 160 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 173
if ($ws.Range('I7').Value2 -ne 529) { $mismatchCount++; Write-Output "MISMATCH Little Village!I7: expected 529 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 542

# East Side
$ws = $wb.Worksheets.Item('East Side')
if ($ws.Range('I6').Value2 -ne 32) { $mismatchCount++; Write-Output "MISMATCH East Side!I6: expected 32 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 33
if ($ws.Range('I7').Value2 -ne 124) { $mismatchCount++; Write-Output "MISMATCH East Side!I7: expected 124 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 125

# Brighton Park
$ws = $wb.Worksheets.Item('Brighton Park')
if ($ws.Range('I6').Value2 -ne 104) { $mismatchCount++; Write-Output "MISMATCH Brighton Park!I6: expected 104 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 105
if ($ws.Range('I7').Value2 -ne 276) { $mismatchCount++; Write-Output "MISMATCH Brighton Park!I7: expected 276 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 277

# North Center
$ws = $wb.Worksheets.Item('North Center')
if ($ws.Range('I6').Value2 -ne 26) { $mismatchCount++; Write-Output "MISMATCH North Center!I6: expected 26 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 27
if ($ws.Range('I7').Value2 -ne 65) { $mismatchCount++; Write-Output "MISMATCH North Center!I7: expected 65 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 66

# Belmont Cragin
$ws = $wb.Worksheets.Item('Belmont Cragin')
if ($ws.Range('I3').Value2 -ne 79) { $mismatchCount++; Write-Output "MISMATCH Belmont Cragin!I3: expected 79 got $($ws.Range('I3').Value2)" }
$ws.Range('I3').Value = 80
if ($ws.Range('I7').Value2 -ne 361) { $mismatchCount++; Write-Output "MISMATCH Belmont Cragin!I7: expected 361 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 362

# West Town
$ws = $wb.Worksheets.Item('West Town')
if ($ws.Range('I2').Value2 -ne 32) { $mismatchCount++; Write-Output "MISMATCH West Town!I2: expected 32 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 33
if ($ws.Range('I7').Value2 -ne 194) { $mismatchCount++; Write-Output "MISMATCH West Town!I7: expected 194 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 195

# United Center
$ws = $wb.Worksheets.Item('United Center')
if ($ws.Range('I2').Value2 -ne 64) { $mismatchCount++; Write-Output "MISMATCH United Center!I2: expected 64 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 65
if ($ws.Range('I7').Value2 -ne 221) { $mismatchCount++; Write-Output "MISMATCH United Center!I7: expected 221 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 222

# Austin
$ws = $wb.Worksheets.Item('Austin')
if ($ws.Range('I2').Value2 -ne 423) { $mismatchCount++; Write-Output "MISMATCH Austin!I2: expected 423 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 424
if ($ws.Range('I6').Value2 -ne 459) { $mismatchCount++; Write-Output "MISMATCH Austin!I6: expected 459 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 460
if ($ws.Range('I7').Value2 -ne 1430) { $mismatchCount++; Write-Output "MISMATCH Austin!I7: expected 1430 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 1432

# Streeterville
$ws = $wb.Worksheets.Item('Streeterville')
if ($ws.Range('I4').Value2 -ne 70) { $mismatchCount++; Write-Output "MISMATCH Streeterville!I4: expected 70 got $($ws.Range('I4').Value2)" }
$ws.Range('I4').Value = 71
if ($ws.Range('I7').Value2 -ne 150) { $mismatchCount++; Write-Output "MISMATCH Streeterville!I7: expected 150 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 151

# North Park
$ws = $wb.Worksheets.Item('North Park')
if ($ws.Range('I2').Value2 -ne 25) { $mismatchCount++; Write-Output "MISMATCH North Park!I2: expected 25 got $($ws.Range('I2').Value2)" }
$ws.Range('I2').Value = 26
if ($ws.Range('I7').Value2 -ne 80) { $mismatchCount++; Write-Output "MISMATCH North Park!I7: expected 80 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 81

# Morgan Park
$ws = $wb.Worksheets.Item('Morgan Park')
if ($ws.Range('I6').Value2 -ne 41) { $mismatchCount++; Write-Output "MISMATCH Morgan Park!I6: expected 41 got $($ws.Range('I6').Value2)" }
$ws.Range('I6').Value = 42
if ($ws.Range('I7').Value2 -ne 134) { $mismatchCount++; Write-Output "MISMATCH Morgan Park!I7: expected 134 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 135

# Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
if ($ws.Range('I5').Value2 -ne 32) { $mismatchCount++; Write-Output "MISMATCH Auburn Gresham!I5: expected 32 got $($ws.Range('I5').Value2)" }
$ws.Range('I5').Value = 33
if ($ws.Range('I7').Value2 -ne 749) { $mismatchCount++; Write-Output "MISMATCH Auburn Gresham!I7: expected 749 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 750

# Beverly
$ws = $wb.Worksheets.Item('Beverly')
if ($ws.Range('I3').Value2 -ne 10) { $mismatchCount++; Write-Output "MISMATCH Beverly!I3: expected 10 got $($ws.Range('I3').Value2)" }
$ws.Range('I3').Value = 11
if ($ws.Range('I7').Value2 -ne 57) { $mismatchCount++; Write-Output "MISMATCH Beverly!I7: expected 57 got $($ws.Range('I7').Value2)" }
$ws.Range('I7').Value = 58

Write-Output "Mismatches: $mismatchCount"
Write-Output "Done"